# Fix: minimal requirements for Render deploy
#
# - Row 3: the phone number in B3 was entered as text; store it as a real
#   number (matches B2's type).
# - Row 4: a new coupon entry (rushika) is appended. The phone number and
#   date must stay as literal text (not auto-converted to a number / date
#   serial), matching how the other rows' Date column is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: convert the stored text "9160057777" into a genuine number.
$ws.Range("B3").Value = 9160057777

# New row 4.
$ws.Range("A4").Value = "rushika"

# B4 must remain literal text "123456789", not a number. Force the cell to
# text format before assigning so Excel doesn't auto-coerce the numeric
# looking string, then restore the default (unstyled) look.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "123456789"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = 5

# D4 must remain literal text "2025-08-16", not a date serial number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2025-08-16"
$ws.Range("D4").Style = "Normal"
